$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E for the new quarters (2018-12-31 and 2018-09-30),
# shifting the existing quarterly data two columns to the right (old D -> new F, etc.)
$ws.Columns("D:E").Insert()

# Copy the number/date formatting from the (now shifted) F:G columns into the
# newly inserted D:E columns so the new cells carry the same styles as the rest
# of the row (date format in header rows, #,##0 style in data rows).
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

function Set-RowValues {
    param([int]$row, [object[]]$vals)
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 4 + $i).Value = $vals[$i]
    }
}

Set-RowValues 7 @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
Set-RowValues 8 @(190400, 150000, 155000, 128300, 140300, 113700, 105000, 99300, 119800, 91700)
Set-RowValues 9 @(152500, 119300, 120400, 101100, 110800, 88300, 82000, 78100, 90600, 71000)
Set-RowValues 10 @(37900, 30700, 34600, 27200, 29500, 25400, 23000, 21200, 29200, 20700)
Set-RowValues 12 @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
Set-RowValues 13 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 14 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 15 @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", 100)
Set-RowValues 17 @(170200, 133300, 134400, 113900, 121100, 97800, 91500, 87800, 101000, 80400)
Set-RowValues 18 @(20200, 16700, 20600, 14400, 19200, 15900, 13500, 11500, 18800, 11300)
Set-RowValues 20 @(500, 3400, 3600, 2200, 1800, 1400, 300, 600, 400, 500)
Set-RowValues 21 @(21900, 20900, 24900, 17000, 21100, 17400, 13900, 12200, 19200, 11900)
Set-RowValues 22 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 23 @(20700, 20100, 24200, 16600, 21000, 17300, 13800, 12100, 19100, 11800)
Set-RowValues 24 @(3800, 4700, 5200, 3400, 6400, 5400, 4400, 3900, 6000, 3700)
Set-RowValues 25 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 26 @(16900, 15400, 19000, 13200, 14600, 11900, 9400, 8200, 13100, 8200)
Set-RowValues 27 @(13400, 12200, 14900, 11200, 10800, 9300, 7700, 6200, 7700, 6200)
Set-RowValues 28 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 29 @(0, "NA", "NA", "NA", -19000, "NA", "NA", "NA", "NA", "NA")
Set-RowValues 30 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 31 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 32 @(-500, -3400, -3600, -2200, -1800, -1400, -300, -600, -400, -500)
Set-RowValues 33 @(13400, 12200, 14900, 11200, -8200, 9300, 7700, 6200, 7700, 6200)
Set-RowValues 34 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 35 @(13400, 12200, 14900, 11200, -8200, 9300, 7700, 6200, 7700, 6200)
Set-RowValues 38 @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
Set-RowValues 41 @(38300, 33100, 37100, 34400, 36700, 20700, 27300, 32300, 35200, 19800)
Set-RowValues 42 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 43 @(4800, 4100, 6200, 1500, 1600, 1800, 2400, 2900, 2400, 2100)
Set-RowValues 44 @(669000, 648200, 581400, 528500, 496100, 478400, 434600, 406200, 410300, 415200)
Set-RowValues 45 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 46 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 47 @(20300, 20800, 19500, 17700, 33800, 15800, "NA", "NA", "NA", "NA")
Set-RowValues 48 @(4700, 4200, 3500, 2100, 800, 700, 700, 800, 900, 900)
Set-RowValues 49 @(1500, 1700, 2000, "NA", 0, "NA", "NA", "NA", "NA", "NA")
Set-RowValues 50 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 51 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 52 @(19900, 32000, 27000, 32200, 34800, 60200, 65000, 69800, 72000, 97000)
Set-RowValues 53 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 54 @(784000, 771000, 705000, 643200, 611000, 606500, 554300, 533500, 541000, 554100)
Set-RowValues 57 @(26100, 35000, 19500, 20200, 22400, 16200, 16300, 13100, 15100, 35700)
Set-RowValues 58 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 59 @(61200, 57500, 55100, 39000, 39900, 39300, 33700, 32000, 28400, 26500)
Set-RowValues 60 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 61 @(208900, 207900, 174300, 144900, 115700, 105200, 84200, 72700, 85900, 91200)
Set-RowValues 62 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 63 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 64 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 65 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 66 @(315700, 315300, 261700, 214800, 194700, 182300, 154400, 141400, 156400, 177500)
Set-RowValues 68 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 69 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 70 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 71 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 72 @(177500, 164200, 152000, 137100, 125900, 134100, 124800, 117100, 110900, 103300)
Set-RowValues 73 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 74 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 75 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 76 @(468400, 455700, 443300, 428400, 416300, 424200, 399900, 392100, 384600, 376600)
Set-RowValues 77 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 80 @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
Set-RowValues 81 @(13400, 12200, 14900, 11200, -8200, 9300, 7700, 6200, 7700, 6200)
Set-RowValues 83 @(1100, 800, 600, 400, 100, 100, 100, 100, 100, 100)
Set-RowValues 84 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 85 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 86 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 87 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 88 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 89 @(-1200, -26000, 8100, -20400, 5200, -26800, -12900, 16600, -1500, 2700)
Set-RowValues 91 @(-1400, -600, -600, -600, -100, "NA", "NA", 0, -200, -400)
Set-RowValues 92 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 93 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 94 @(-1400, -2100, -26700, -600, -200, -300, 0, 0, -100, -200)
Set-RowValues 96 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 97 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 98 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 99 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 100 @(-700, 31900, 21000, 19500, 10600, 20900, 6800, -18300, -5900, 16400)
Set-RowValues 101 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 102 @(-3300, 3800, 2400, -1400, 15200, -6700, -6100, -1800, -7500, 18900)

Write-Host "Done applying GRBK quarterly financial updates"
